# dv_power BOM update (0.1.14): added SMD connector (X1, X3 / PLD-20),
# PCB redrawn -> several parts renumbered/renamed, R-rows merged into one.
#
# Note: writing .Value on this sheet's text cells (style 3, quotePrefix
# text style) resets the cell style to the plain bordered style (4), so
# after updating values we restore the correct per-column style via
# Copy/PasteSpecial(formats) from untouched scratch cells captured up
# front (G1 = style 3 "text" look, G2 = style 4 "numeric/blank" look).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture pristine formats into scratch cells before any edits ---
$ws.Range("A2").Copy()
$ws.Range("G1").PasteSpecial(-4122)   # style 3: bordered text cell
$ws.Range("D2").Copy()
$ws.Range("G2").PasteSpecial(-4122)   # style 4: bordered blank/number cell

# --- Row 2 (X2) : unchanged content, Value (E) stays blank ---
$ws.Range("A2").Value = "X2"
$ws.Range("B2").Value = "Header, 20-Pin, Dual row"
$ws.Range("C2").Value = "PLD2-20"
$ws.Range("D2").Value = 1

# --- Row 3 (new): X1, X3 / Header, 20-Pin, Dual row / PLD-20 ---
$ws.Range("A3").Value = "X1, X3"
$ws.Range("B3").Value = "Header, 20-Pin, Dual row"
$ws.Range("C3").Value = "PLD-20"
$ws.Range("D3").Value = 2

# --- Row 4: VT1, VT2 (was row 3) ---
$ws.Range("A4").Value = "VT1, VT2"
$ws.Range("B4").Value = "Биполярный транзистор"
$ws.Range("C4").Value = '2Т3117А/ПК "ОСМ"'
$ws.Range("D4").Value = 2

# --- Row 5: VD1, VD2 (was row 4) ---
$ws.Range("A5").Value = "VD1, VD2"
$ws.Range("B5").Value = "Сдвоенный диод"
$ws.Range("C5").Value = "2Д222ВС ОСМ"
$ws.Range("D5").Value = 2

# --- Row 6: U1, U2 (was row 5), now has a Наименование value, Value blank ---
$ws.Range("A6").Value = "U1, U2"
$ws.Range("B6").Value = "DC-DC преобразователь"
$ws.Range("C6").Value = "СМПВ 1.5 5.0 ОВ"
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = ""

# --- Row 7: R1, R2, R3, R4 (merged former rows 6 & 7) ---
$ws.Range("A7").Value = "R1, R2, R3, R4"
$ws.Range("B7").Value = "SMD-резистор"
$ws.Range("C7").Value = "ОСМ Р1-12"
$ws.Range("D7").Value = 4
$ws.Range("E7").Value = "1к"

# --- Row 8: FU1..FU4 (unchanged) ---
$ws.Range("A8").Value = "FU1, FU2, FU3, FU4"
$ws.Range("B8").Value = "Предохранитель плавкий"
$ws.Range("C8").Value = 'ВП1-2 "ВП"'
$ws.Range("D8").Value = 4
$ws.Range("E8").Value = "0.25A"

# --- Row 9: C2, C4, C6, C8 - Наименование text updated ---
$ws.Range("A9").Value = "C2, C4, C6, C8"
$ws.Range("B9").Value = "Поляризованный SMD-конденсатор"
$ws.Range("C9").Value = "ОС К53-68 "
$ws.Range("D9").Value = 4
$ws.Range("E9").Value = "10мкФх50В"

# --- Row 10: C1, C3, C5, C7 - Наименование text updated ---
$ws.Range("A10").Value = "C1, C3, C5, C7"
$ws.Range("B10").Value = "SMD-конденсатор"
$ws.Range("C10").Value = "К10-79"
$ws.Range("D10").Value = 4
$ws.Range("E10").Value = "0,1мкФ"

# --- restore the text style (3) on every A/B/C cell and the text E cells ---
$ws.Range("G1").Copy()
$ws.Range("A2:C10").PasteSpecial(-4122)
$ws.Range("E7:E10").PasteSpecial(-4122)

# --- restore the blank/numeric style (4) on D column and the blank E cells ---
$ws.Range("G2").Copy()
$ws.Range("D2:D10").PasteSpecial(-4122)
$ws.Range("E2:E6").PasteSpecial(-4122)

# --- clean up scratch cells ---
$ws.Range("G1:G2").Clear()
